$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.332.04"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.178.49"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'252.55"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").Value = "'0.605"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").Value = "'73.28"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.582"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").Value = "'40.13"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "'6.73"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Value = "2.504.92"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "'14.17"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").Value = "2.167.57"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").Value = "'0.769"
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").Value = "42.237.22"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "'70.61"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "'5.84"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").Value = "'226.84"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "'9.32"
$ws.Range("E23").Value = "  -4.58%  "
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "'10.47"
$ws.Range("E26").Value = "  -4.22%  "
$ws.Range("D27").Value = "'3.38"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.18"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.16"
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("D30").Value = "'36.61"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "'169.85"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D32").Value = "'19.99"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "'0.0817"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").Value = "'5.11"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").Value = "'4.19"
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("D38").Value = "'0.0335"
$ws.Range("E38").Value = "  +3.48%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'11.76"
$ws.Range("E39").Value = "  -8.26%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.05"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.195"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'59.15"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").Value = "'5.14"
$ws.Range("E43").Value = "  -6.84%  "
$ws.Range("D44").Value = "'101.83"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("E45").Value = "  +6.66%  "
$ws.Range("D46").Value = "'0.0970"
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("B47").Value = "WOONetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D47").Value = "'0.458"
$ws.Range("E47").Value = "  +5.06%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'8.14"
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("E50").Value = "  -1.35%  "
